$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.741320133209229
$ws.Range("B1").Value = 2.300427436828613
$ws.Range("C1").Value = 2.234244346618652
$ws.Range("D1").Value = 6.050808906555176
$ws.Range("E1").Value = 0.7055863738059998
